# Generate Report for Handback
# Adds a new row (for file 597be984-7d85-4f6b-823c-5c225f7a940d.md) to the
# "Overview", "zh-cn" and "de-de" tables/sheets, growing each table from
# A1:*3 to A1:*4.

$wb = $excel.ActiveWorkbook

$srcFile      = "597be984-7d85-4f6b-823c-5c225f7a940d.md"
$srcFileE2e   = "e2e\" + $srcFile
$statusInSync = "Handed back: in sync with en-US"
$ext          = ".md"
$sourcePath   = "e2e"
$priority     = "ht"
$trueStr      = "True"
$falseStr     = "False"

$zhXlf        = "597be984-7d85-4f6b-823c-5c225f7a940d.fa31be5988ecab6d5574378420afa35ce517060f.zh-cn.xlf"
$zhHandoffDt  = "2016-08-19 00:42:59"
$zhHandbackDt = "2016-08-19 00:43:27"

$deXlf        = "597be984-7d85-4f6b-823c-5c225f7a940d.fa31be5988ecab6d5574378420afa35ce517060f.de-de.xlf"
$deHandoffDt  = "2016-08-19 00:43:08"
$deHandbackDt = "2016-08-19 00:43:34"

$srcRepoUrl   = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6a0b1fbc6a5c6c6f5b6c1e2f3a4b5c6d7e8f9a0b/e2e/" + $srcFile
$zhRepoUrl    = "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/1a2b3c4d5e6f7a8b9c0d1e2f3a4b5c6d7e8f9a0b/e2e/" + $srcFile
$deRepoUrl    = "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/0b1c2d3e4f5a6b7c8d9e0f1a2b3c4d5e6f7a8b9c/e2e/" + $srcFile

# ---------------------------------------------------------------------
# Overview sheet (table3 -> "Overview", columns A:G)
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.ListRows.Add() | Out-Null

$wsOverview.Range("A4").Value = $srcFile
$wsOverview.Range("B4").Value = $srcFileE2e
$wsOverview.Range("C4").Value = $ext
$wsOverview.Range("E4").Value = $statusInSync
$wsOverview.Range("F4").Value = $statusInSync
$wsOverview.Range("G4").Value = $deHandoffDt

$wsOverview.Hyperlinks.Add($wsOverview.Range("B4"), $srcRepoUrl, "", "", $srcFileE2e) | Out-Null

# ---------------------------------------------------------------------
# zh-cn sheet (table1, columns A:P)
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$loZh = $wsZh.ListObjects.Item(1)
$loZh.ListRows.Add() | Out-Null

$wsZh.Range("A4").Value = $srcFile
$wsZh.Range("B4").Value = $ext
$wsZh.Range("C4").Value = $statusInSync
$wsZh.Range("D4").Value = $sourcePath
$wsZh.Range("E4").Value = $priority
$wsZh.Range("F4").Value = $trueStr
$wsZh.Range("G4").Value = $zhXlf
$wsZh.Range("H4").Value = $zhHandoffDt
$wsZh.Range("I4").Value = $srcFile
$wsZh.Range("J4").Value = $zhXlf
$wsZh.Range("K4").Value = $zhHandbackDt
$wsZh.Range("L4").Value = ""
$wsZh.Range("M4").Value = $trueStr
$wsZh.Range("N4").Value = ""
$wsZh.Range("O4").Value = $falseStr
$wsZh.Range("P4").Value = ""

$wsZh.Hyperlinks.Add($wsZh.Range("A4"), $srcRepoUrl, "", "", $srcFile) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("I4"), $zhRepoUrl, "", "", $srcFile) | Out-Null

# ---------------------------------------------------------------------
# de-de sheet (table2, columns A:P)
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$loDe = $wsDe.ListObjects.Item(1)
$loDe.ListRows.Add() | Out-Null

$wsDe.Range("A4").Value = $srcFile
$wsDe.Range("B4").Value = $ext
$wsDe.Range("C4").Value = $statusInSync
$wsDe.Range("D4").Value = $sourcePath
$wsDe.Range("E4").Value = $priority
$wsDe.Range("F4").Value = $trueStr
$wsDe.Range("G4").Value = $deXlf
$wsDe.Range("H4").Value = $deHandoffDt
$wsDe.Range("I4").Value = $srcFile
$wsDe.Range("J4").Value = $deXlf
$wsDe.Range("K4").Value = $deHandbackDt
$wsDe.Range("L4").Value = ""
$wsDe.Range("M4").Value = $trueStr
$wsDe.Range("N4").Value = ""
$wsDe.Range("O4").Value = $falseStr
$wsDe.Range("P4").Value = ""

$wsDe.Hyperlinks.Add($wsDe.Range("A4"), $srcRepoUrl, "", "", $srcFile) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("I4"), $deRepoUrl, "", "", $srcFile) | Out-Null
